$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5201.25
$ws.Range("I19").Value = 9812.5
$ws.Range("K19").Value = 9812.5
$ws.Range("M19").Value = -9637.5

$ws.Range("H113").Value = 3544
$ws.Range("I113").Value = 3050
$ws.Range("J113").Value = 4136.8
$ws.Range("K113").Value = 3050
$ws.Range("L113").Value = 4136.8
$ws.Range("M113").Value = 204
$ws.Range("N113").Value = -10644.8

$ws.Range("H116").Value = 3543.6428
$ws.Range("I116").Value = 1233
$ws.Range("J116").Value = 4173.8184
$ws.Range("K116").Value = 1233
$ws.Range("L116").Value = 4173.8184
$ws.Range("M116").Value = 2209
$ws.Range("N116").Value = -11057.8184

$ws.Range("H129").Value = 189479.8
$ws.Range("I129").Value = 339.8
$ws.Range("J129").Value = 209181.88
$ws.Range("K129").Value = 1019.4
$ws.Range("L129").Value = 627545.64
$ws.Range("M129").Value = 3980.6
$ws.Range("N129").Value = -637545.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2020.7826
$ws.Range("I122").Value = 2031.619
$ws.Range("J122").Value = 1907
$ws.Range("K122").Value = 6094.857
$ws.Range("L122").Value = 5721
$ws.Range("M122").Value = -3644.857
$ws.Range("N122").Value = -10621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10451.484
$ws.Range("I31").Value = 12110.667
$ws.Range("K31").Value = 12110.667
$ws.Range("M31").Value = -11815.667

$ws.Range("H34").Value = 10451.484
$ws.Range("I34").Value = 12110.667
$ws.Range("K34").Value = 12110.667
$ws.Range("M34").Value = -11908.667

$ws.Range("H94").Value = 7794.923
$ws.Range("I94").Value = 2400
$ws.Range("J94").Value = 12419.143
$ws.Range("K94").Value = 2400
$ws.Range("L94").Value = 12419.143
$ws.Range("M94").Value = -1949
$ws.Range("N94").Value = -13321.143

$ws.Range("H99").Value = 15155478
$ws.Range("I99").Value = 3212.818
$ws.Range("K99").Value = 3212.818
$ws.Range("M99").Value = -1714.818

$ws.Range("H122").Value = 1594.6875
$ws.Range("I122").Value = 1901.5
$ws.Range("J122").Value = 1083.3334
$ws.Range("K122").Value = 5704.5
$ws.Range("L122").Value = 3250.0002
$ws.Range("M122").Value = -3254.5
$ws.Range("N122").Value = -8150.0002

$ws.Range("H126").Value = 15155478
$ws.Range("I126").Value = 3212.818
$ws.Range("K126").Value = 9638.454000000002
$ws.Range("M126").Value = -7168.454000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 963.1667
$ws.Range("J11").Value = 1866.6666
$ws.Range("L11").Value = 5599.9998
$ws.Range("N11").Value = -5879.9998

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H75").Value = 1748.7273
$ws.Range("I75").Value = 1198.2
$ws.Range("J75").Value = 2207.5
$ws.Range("K75").Value = 3594.6
$ws.Range("L75").Value = 6622.5
$ws.Range("M75").Value = -2596.6
$ws.Range("N75").Value = -8618.5

$ws.Range("H78").Value = 1748.7273
$ws.Range("I78").Value = 1198.2
$ws.Range("J78").Value = 2207.5
$ws.Range("K78").Value = 10783.8
$ws.Range("L78").Value = 19867.5
$ws.Range("M78").Value = -5791.800000000001
$ws.Range("N78").Value = -29851.5

$ws.Range("H114").Value = 1005
$ws.Range("I114").Value = 1126.2222
$ws.Range("J114").Value = 868.625
$ws.Range("K114").Value = 3378.6666
$ws.Range("L114").Value = 2605.875
$ws.Range("M114").Value = -124.6665999999996
$ws.Range("N114").Value = -9113.875

$ws.Range("H116").Value = 774.75
$ws.Range("I116").Value = 774.75
$ws.Range("K116").Value = 2324.25
$ws.Range("M116").Value = 1117.75

$ws.Range("H117").Value = 2337.9375
$ws.Range("I117").Value = 889.25
$ws.Range("J117").Value = 2820.8333
$ws.Range("K117").Value = 2667.75
$ws.Range("L117").Value = 8462.499899999999
$ws.Range("M117").Value = 774.25
$ws.Range("N117").Value = -15346.4999

$ws.Range("H120").Value = 9912.857
$ws.Range("I120").Value = 5878
$ws.Range("K120").Value = 17634
$ws.Range("M120").Value = -12796

$ws.Range("H129").Value = 500819
$ws.Range("J129").Value = 625849.9
$ws.Range("L129").Value = 1877549.7
$ws.Range("N129").Value = -1887549.7

$ws.Range("H137").Value = 2729.6
$ws.Range("J137").Value = 4493.2
$ws.Range("L137").Value = 13479.6
$ws.Range("N137").Value = -23679.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11691.308
$ws.Range("I70").Value = 17605.285
$ws.Range("K70").Value = 17605.285
$ws.Range("M70").Value = -17335.285

$ws.Range("H73").Value = 11691.308
$ws.Range("I73").Value = 17605.285
$ws.Range("K73").Value = 17605.285
$ws.Range("M73").Value = -16669.285

$ws.Range("H102").Value = 22729472
$ws.Range("J102").Value = 976.8
$ws.Range("L102").Value = 976.8
$ws.Range("N102").Value = -4220.8

$ws.Range("H113").Value = 2845.3333
$ws.Range("I113").Value = 2131.111
$ws.Range("J113").Value = 3916.6667
$ws.Range("K113").Value = 2131.111
$ws.Range("L113").Value = 3916.6667
$ws.Range("M113").Value = 38.88900000000012
$ws.Range("N113").Value = -8256.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6476.4707
$ws.Range("I7").Value = 7030.385
$ws.Range("J7").Value = 4676.25
$ws.Range("K7").Value = 7030.385
$ws.Range("L7").Value = 4676.25
$ws.Range("M7").Value = -6918.385
$ws.Range("N7").Value = -4900.25

$ws.Range("H40").Value = 4298
$ws.Range("I40").Value = 3291.6
$ws.Range("J40").Value = 4633.467
$ws.Range("K40").Value = 3291.6
$ws.Range("L40").Value = 4633.467
$ws.Range("M40").Value = -3155.6
$ws.Range("N40").Value = -4905.467

$ws.Range("H122").Value = 1035076.7
$ws.Range("I122").Value = 1403010.9
$ws.Range("J122").Value = 4861
$ws.Range("K122").Value = 4209032.699999999
$ws.Range("L122").Value = 14583
$ws.Range("M122").Value = -4206582.699999999
$ws.Range("N122").Value = -19483

$ws.Range("H126").Value = 6476.4707
$ws.Range("I126").Value = 7030.385
$ws.Range("J126").Value = 4676.25
$ws.Range("K126").Value = 21091.155
$ws.Range("L126").Value = 14028.75
$ws.Range("M126").Value = -18621.155
$ws.Range("N126").Value = -18968.75

$ws.Range("H132").Value = 2169.074
$ws.Range("I132").Value = 1474.5294
$ws.Range("J132").Value = 3349.8
$ws.Range("K132").Value = 4423.5882
$ws.Range("L132").Value = 10049.4
$ws.Range("M132").Value = -1893.5882
$ws.Range("N132").Value = -15109.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1528.8889
$ws.Range("I122").Value = 1589.24
$ws.Range("K122").Value = 4767.72
$ws.Range("M122").Value = -2317.72

$ws.Range("H126").Value = 1482.3334
$ws.Range("I126").Value = 1997.375
$ws.Range("J126").Value = 1070.3
$ws.Range("K126").Value = 5992.125
$ws.Range("L126").Value = 3210.9
$ws.Range("M126").Value = -3522.125
$ws.Range("N126").Value = -8150.9

$ws.Range("H127").Value = 25107.25
$ws.Range("J127").Value = 25107.25
$ws.Range("L127").Value = 25107.25
$ws.Range("N127").Value = -35027.25

$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
